$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking values
# (e.g. "299.85") are stored as text, matching the inline-string
# source data, rather than being auto-converted to numbers.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "42.232.13"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "2.272.12"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "299.85"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "95.60"
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("D7").Value = "0.496"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").Value = "33.14"
$ws.Range("E10").Value = "  -4.26%  "
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -6.60%  "
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").Value = "15.93"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "2.625.47"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "2.290.42"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "0.787"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "42.174.79"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").Value = "11.69"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "66.38"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").Value = "235.22"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("E28").Value = "  -4.59%  "
$ws.Range("D29").Value = "167.21"
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("E30").Value = "  -4.69%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "33.54"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "9.10"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("E34").Value = "  +6.75%  "
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "16.70"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "2.79"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").Value = "0.0988"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("E42").Value = "  -4.08%  "
$ws.Range("E43").Value = "  -7.19%  "
$ws.Range("D44").Value = "1.957.74"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "0.0278"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "17.49"
$ws.Range("E46").Value = "  -5.09%  "
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("D49").Value = "2.496.79"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").Value = "52.13"
$ws.Range("E50").Value = "  -5.52%  "
$ws.Range("E51").Value = "  -3.25%  "

# Restore the original General number format / default style so the
# only observable change is the cell text, not formatting metadata.
$colD.NumberFormat = "General"
$colD.Style = "Normal"
